# Status and query sheet update for syncup
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Status(Summary)")

# New weekly status row (row 19) appended under Table6 (A1:D18 -> A1:D19)
$row = 19

$ws.Cells.Item($row, 1).Value = 43986
$ws.Cells.Item($row, 2).Value = "Implementation"
$ws.Cells.Item($row, 3).Value = "App design and study"

$comment = "1) Advised to use Influenzanet app directly rather than new app+integration`n2) Completed some of classification of group components related to roles and also question fetch from flat rendered items`n3) Started migration of all the app tointegrate with the Influenzanet app`n4) Every question should not be displayed in one screen i.e a no to question per screen. Needs to be clarified in the meeting"

$ws.Cells.Item($row, 4).Value = $comment

# Bold the two emphasised words ("not" / "no") within the comment, matching
# the rest of the run's existing (non-bold) Calibri 11 formatting
$cell = $ws.Cells.Item($row, 4)
$cell.Characters(298, 4).Font.Bold = $true
$cell.Characters(335, 3).Font.Bold = $true

# Copy the formatting (number format / alignment / wrap) from the row above,
# which already carries the styles used throughout the table body
$ws.Range("A18:D18").Copy()
$ws.Range("A19:D19").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Row auto-sizes to the wrapped comment text, same as the rest of the table
$ws.Rows.Item($row).RowHeight = 86

# Grow Table6 so the new row participates in the table/autofilter range
$table = $ws.ListObjects.Item("Table6")
$table.Resize($ws.Range("A1:D19"))

# Update selection to mirror the authored workbook view
$ws.Range("D19").Select()
